$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54 data
$ws.Cells.Item(54, 1).Value = 45986
$ws.Cells.Item(54, 2).Value = 2025
$ws.Cells.Item(54, 3).Value = 2.622852459381209
$ws.Cells.Item(54, 4).Value = 2026
$ws.Cells.Item(54, 5).Value = 1.946625946175717

# Copy style from the cell above (A53) so A54 keeps the date style (s="2")
$ws.Cells.Item(53, 1).Copy() | Out-Null
$ws.Cells.Item(54, 1).PasteSpecial(-4122) | Out-Null # xlPasteFormats
